$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move the existing headers over one column: A1/B1/C1 -> B1/C1/D1
$ws.Range("A1").Value = $null
$ws.Range("B1").Value = "fuel change"
$ws.Range("C1").Value = "velocity change"
$ws.Range("D1").Value = "multiplier"

# New rows
$ws.Range("A2").Value = "click"
$ws.Range("B2").Value = 0.25
$ws.Range("C2").Value = 2
$ws.Range("D2").Formula = "=C2/B2"

$ws.Range("A3").Value = "uploop/s"
$ws.Range("B3").Value = 4.8
$ws.Range("C3").Value = 38.4
$ws.Range("D3").Formula = "=C3/B3"

$ws.Range("I21").Select()
